$wb = $excel.ActiveWorkbook

$sheetNames = @("2025", "2030", "2035", "2040", "2045", "2050")

$values = @(
    @{ A = 4830.603421671776; B = 6914.260623028063; E = 15339.30996716316; I = 36843.73937326; M = 11716.236729175; N = 4270.376965689151; O = 6924.712879078099 },
    @{ A = 4830.603421671776; B = 8041.469130500362; E = 41322.54231073364; I = 67177.40575988838; M = 22588.905778934; N = 10465.8025602956; O = 12108.14566538872 },
    @{ A = 4830.603421671776; B = 8041.469130500362; E = 61689.09165384799; G = 7864.06113287106; I = 85324.99695076495; L = 4893.458355751803; M = 28513.42754631791; N = 13431.99981023839; O = 15331.82941245748 },
    @{ A = 4830.603421671776; B = 8041.469130500362; E = 61689.09165384799; G = 7864.06113287106; I = 85324.99695076495; L = 4893.458355751803; M = 28513.42754631791; N = 13431.99981023839; O = 15331.82941245748 },
    @{ A = 4830.603421671776; B = 8041.469130500362; E = 61689.09165384799; G = 7864.06113287106; I = 85324.99695076495; L = 4893.458355751803; M = 28513.42754631791; N = 13431.99981023839; O = 15331.82941245748 },
    @{ A = 4830.603421671776; B = 8041.469130500362; E = 61689.09165384799; G = 7864.06113287106; I = 85324.99695076495; L = 4893.458355751803; M = 28513.42754631791; N = 13431.99981023839; O = 15331.82941245748 }
)

for ($i = 0; $i -lt $sheetNames.Length; $i++) {
    $ws = $wb.Worksheets.Item([string]$sheetNames[$i])
    $cols = $values[$i]
    foreach ($col in $cols.Keys) {
        $ws.Range($col + "2").Value = $cols[$col]
    }
}
